$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force literal text so numeric-looking strings (e.g. "208.96")
    # aren't auto-coerced into Excel numbers, then restore the
    # original (unstyled) cell format so no stray style is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.413.03"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.564.98"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue $ws.Range("D5") "208.96"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.21%  "
Set-TextValue $ws.Range("D8") "21.99"
$ws.Range("E8").Value = "  -1.29%  "
Set-TextValue $ws.Range("D10") "0.0591"
$ws.Range("E10").Value = "  -0.12%  "
Set-TextValue $ws.Range("D11") "0.0866"
$ws.Range("D12").Value = "1.787.66"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Value = "1.550.43"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  -2.71%  "
Set-TextValue $ws.Range("D16") "63.47"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "27.413.22"
Set-TextValue $ws.Range("D18") "212.90"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("E19").Value = "  -0.89%  "
Set-TextValue $ws.Range("D20") "7.25"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -0.64%  "
Set-TextValue $ws.Range("D23") "9.51"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +2.11%  "
Set-TextValue $ws.Range("D25") "152.98"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.22%  "
Set-TextValue $ws.Range("D27") "6.71"
$ws.Range("E27").Value = "  -0.33%  "
Set-TextValue $ws.Range("D28") "14.97"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +0.77%  "
Set-TextValue $ws.Range("D32") "3.20"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "1.370.33"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +1.38%  "
Set-TextValue $ws.Range("D36") "0.963"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  +0.96%  "
Set-TextValue $ws.Range("D39") "0.531"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  -0.19%  "
Set-TextValue $ws.Range("D42") "0.974"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +0.59%  "
Set-TextValue $ws.Range("D44") "63.92"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "1.700.28"
$ws.Range("E47").Value = "  -1.22%  "
Set-TextValue $ws.Range("D48") "85.49"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").Value = "0.0₇0991"
$ws.Range("E49").Value = "  -2.36%  "
Set-TextValue $ws.Range("D50") "0.0954"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  -0.81%  "
